$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($old, $new) {
  $rng = $d.Content
  $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null
}

function Delete-ParagraphContaining($searchText) {
  $rng = $d.Content
  $found = $rng.Find.Execute($searchText)
  if ($found) {
    $rng.Expand(4) | Out-Null
    $rng.Delete() | Out-Null
  }
}

# ---------------------------------------------------------------------
# 1. Update the git revision string
# ---------------------------------------------------------------------
Replace-Text "a4b3d5b" "a291098"

# ---------------------------------------------------------------------
# 2. Reword the Methods paragraph, pointing the endogeneity aside at the
#    new Protocol Deviations section
# ---------------------------------------------------------------------
$oldMethods = "Except as noted below, we analyzed the data as specified in the protocol using Stata 16 (StataCorp LLC, College Station, Texas, USA). Briefly, we analyzed resource use (person-hours) on the log scale using extended interval regression (eintreg) and used a likelihood-adjusted-censoring inverse-probability-weighted regression adjustment model (LAC-IPWRA; stteffects) to estimate mean differences in time-to-completion. Ongoing reviews were right censored at the end of data collection (31 January 2023) and all analyses accounted for this censoring. We had no reason to suspect informative (nonrandom) censoring, so did not model a censoring mechanism. Because we did not randomize reviews to use recommended ML versus no ML (for example), we modelled ML use as an endogenously assigned treatment predicted by field (healthcare or welfare) and pre-specification (existence of a protocol), as planned, in all analyses except that for the secondary analysis recommended versus non-recommended ML use with respect to resource use. While there was some statistically significant evidence of endogeneity from the corresponding time-to-completion analysis and an exploratory logistic regression, the estimate of relative resource use obtained using the planned model appeared to dramatically overestimate the effect of recommended ML use. We therefore used a model for this analysis that did not account for possible endogeneity. We re-expressed all estimates as ratios (relative resource use and relative time-to-completion) to aid generalization to other institutions. We present two-sided 95% confidence intervals and p-values where appropriate and use a prespecified p < 0.05 significance criterion throughout. We also present the time-to-completion data using Kaplan-Meier estimates of survivor functions."
$newMethods = "Except as noted, all statistical analyses were performed as specified in our protocol using Stata 16 (StataCorp LLC, College Station, Texas, USA). Briefly, we analyzed resource use (person-hours) on the log scale using extended interval regression (eintreg) and used a likelihood-adjusted-censoring inverse-probability-weighted regression adjustment model (LAC-IPWRA; stteffects) to estimate mean differences in time-to-completion. Ongoing reviews were right censored at the end of data collection (31 January 2023) and all analyses accounted for this censoring. We had no reason to suspect informative (nonrandom) censoring, so did not model a censoring mechanism. Because we did not randomize reviews to use recommended ML versus no ML (for example), we modelled ML use as an endogenously assigned treatment predicted by field (healthcare or welfare) and pre-specification (existence of a protocol), as planned, in all but one analysis (see Protocol Deviations). We re-expressed all estimates as ratios (relative resource use and relative time-to-completion) to aid generalization to other institutions. We present two-sided 95% confidence intervals and p-values where appropriate and use a prespecified p < 0.05 significance criterion throughout. We also present the time-to-completion data using Kaplan-Meier estimates of survivor functions."
Replace-Text $oldMethods $newMethods

# ---------------------------------------------------------------------
# 3. Insert the new 'Protocol Deviations' section right after Methods
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("time-to-completion data using Kaplan-Meier estimates of survivor functions.") | Out-Null
$methodsPara = $anchor.Paragraphs(1)

$methodsPara.Range.InsertParagraphAfter() | Out-Null
$headingPara = $methodsPara.Next()
$headingPara.Range.Text = "Protocol Deviations"
$headingPara.Style = "Heading1"

$headingPara.Range.InsertParagraphAfter() | Out-Null
$body1Para = $headingPara.Next()
$body1Para.Style = "Normal"
$body1Para.Range.Text = "We had planned to model ML use as an endogenously assigned treatment in all analyses. However, we chose to deviate from protocol for the secondary analysis of recommended versus non-recommended ML use for the outcome of resource use. While there was some statistically significant evidence of endogeneity from the corresponding time-to-completion analysis and an exploratory logistic regression, the estimate of relative resource use obtained using the planned model appeared to dramatically overestimate the effect of recommended ML use. We therefore used a model for this analysis that did not account for possible endogeneity."
$body1Para.Alignment = 3

$body1Para.Range.InsertParagraphAfter() | Out-Null
$body2Para = $body1Para.Next()
$body2Para.Style = "Normal"
$body2Para.Range.Text = "We updated the preprint version of the protocol during data extraction but before starting the analysis or unblinding the statistician (CJR) to redefine the comparisons in terms of under- and overuse of machine learning. However, only two reviews were judged to have under- or overused machine learning, so it was not possible to perform the revised analyses. We therefore performed and report the analyses as originally planned."
$body2Para.Alignment = 3

# ---------------------------------------------------------------------
# 4. Round the summary-table ratio estimates to 2 s.f. / 1 d.p.
# ---------------------------------------------------------------------
Replace-Text "0.10 (0.06 to 0.17)" "0.1 (0.1 to 0.2)"
Replace-Text "1.49 (0.78 to 2.85)" "1.5 (0.8 to 2.9)"
Replace-Text "0.08 (0.04 to 0.20)" "0.1 (0.0 to 0.2)"
Replace-Text "1.10 (0.69 to 1.78)" "1.1 (0.7 to 1.8)"
Replace-Text "0.89 (0.55 to 1.42)" "0.9 (0.6 to 1.4)"
Replace-Text "1.10 (0.68 to 1.79)" "1.1 (0.7 to 1.8)"

# Table footnote
$oldFootnote = "¹Data are means of samples restricted to completed reviews and may underestimate resource use (person-hours) and time-to-completion (weeks) due to right-censoring of ongoing projects. ²Estimates are relative resource use and relative time-to-completion, account for right-censored outcomes and nonrandom endogenous treatment allocation, and are adjusted for planned meta-analysis."
$newFootnote = "¹Data are means of samples restricted to completed (uncensored) reviews. ²Estimates are relative resource use and relative time-to-completion, account for right-censored outcomes and, except for the recommended versus non-recommended ML use comparison for the outcome resource use, also account for nonrandom endogenous treatment allocation. All estimates are adjusted for planned meta-analysis."
Replace-Text $oldFootnote $newFootnote

# ---------------------------------------------------------------------
# 5. Drop the old References / Appendix-1 / Appendix-2 scaffolding; the
#    Protocol Deviations content now lives in the Methods section and
#    the old Appendix 2 becomes the sole appendix.
# ---------------------------------------------------------------------
Delete-ParagraphContaining "TODO: Add references."
Delete-ParagraphContaining "Appendix 1 — Protocol Deviations"
Delete-ParagraphContaining "overuse of machine learning (TODO: Cite revision)"
Delete-ParagraphContaining "Appendix 2 — Full Regression Results"

Replace-Text "References" "Appendix — Full Regression Results"

